$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 171 - this shifts rows 171:222 down to 172:223
# and pushes their formatting/values along with them, matching the
# target diff where every row from 171 onward is shifted down by one
# and a brand new record is inserted at the top of that shift.
$ws.Rows("171:171").Insert()

# Populate the newly inserted row 171 with the new weekly record.
$ws.Range("A171").Value = 6
$ws.Range("B171").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C171").Value = "Metropolitana"
$ws.Range("D171").Value = 44524
$ws.Range("D171").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E171").Value = 13
$ws.Range("F171").Value = "Fruta"
$ws.Range("G171").Value = 100101
$ws.Range("H171").Value = "Berries"
$ws.Range("I171").Value = 100101001
$ws.Range("J171").Value = "Arándano (blue)"
$ws.Range("K171").Value = "Sin especificar"
$ws.Range("L171").Value = "Primera"
$ws.Range("M171").Value = 2500
$ws.Range("N171").Value = 5000
$ws.Range("O171").Value = 5000
$ws.Range("P171").Value = 5000
$ws.Range("Q171").Value = "$/bandeja 2 kilos"
$ws.Range("R171").Value = "Provincia de Curicó"
$ws.Range("S171").Value = 2500
$ws.Range("T171").Value = 2
